$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.344.78"
$ws.Range("E2").Value = "  -4.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.977.07"
$ws.Range("E3").Value = "  -5.80%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.21"
$ws.Range("E5").Value = "  -2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  -4.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.11"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.79"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("E10").Value = "  -4.77%  "
$ws.Range("E11").Value = "  -6.31%  "
$ws.Range("E12").Value = "  -6.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.877"
$ws.Range("E13").Value = "  -5.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.02"
$ws.Range("E14").Value = "  -7.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.268.14"
$ws.Range("E15").Value = "  -5.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.14"
$ws.Range("E16").Value = "  -7.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.969.57"
$ws.Range("E17").Value = "  -6.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.283.30"
$ws.Range("E18").Value = "  -4.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.72"
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.28"
$ws.Range("E20").Value = "  -4.91%  "
$ws.Range("E21").Value = "  -6.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.03"
$ws.Range("E22").Value = "  -4.03%  "
$ws.Range("E23").Value = "  -9.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  -4.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  +3.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.15"
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.98"
$ws.Range("E28").Value = "  -7.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.20"
$ws.Range("E29").Value = "  -8.61%  "
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.11"
$ws.Range("E31").Value = "  -5.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.73"
$ws.Range("E32").Value = "  -10.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0580"
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("E34").Value = "  +7.59%  "
$ws.Range("E35").Value = "  -11.48%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("E38").Value = "  -11.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.79"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("E40").Value = "  -8.71%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0205"
$ws.Range("E42").Value = "  -6.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.06"
$ws.Range("E43").Value = "  -8.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.369.91"
$ws.Range("E44").Value = "  -3.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0871"
$ws.Range("E45").Value = "  -9.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.68"
$ws.Range("E46").Value = "  -8.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.29"
$ws.Range("E47").Value = "  -5.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.20"
$ws.Range("E48").Value = "  -5.96%  "
$ws.Range("E50").Value = "  -9.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.27"
$ws.Range("E51").Value = "  -3.95%  "
